$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Maybe some cool analysis project"
$ws.Range("B16").Select()
